$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on every cell we touch so numeric-looking strings
# (prices like '5.02', '1.00', '0.0000188') are preserved exactly as
# literal text instead of being auto-coerced into Doubles (which would
# introduce float rounding noise or scientific notation).
$targetCells = @("D2", "E2", "D3", "E3", "E4", "D5", "E5", "D6", "E6", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "E21", "D22", "E22", "D23", "E23", "B24", "C24", "D24", "E24", "B25", "C25", "D25", "E25", "E26", "D27", "E27", "D28", "E28", "D29", "E29", "E30", "D31", "E31", "D32", "E32", "D33", "E33", "E34", "D35", "E35", "E36", "D37", "E37", "D38", "E38", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "D48", "E48", "E50", "E51")
foreach ($ref in $targetCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "71.503.51"
$ws.Range("E2").Value = "  +4.10%  "
$ws.Range("D3").Value = "2.627.37"
$ws.Range("E3").Value = "  +3.99%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "607.07"
$ws.Range("E5").Value = "  +2.13%  "
$ws.Range("D6").Value = "181.52"
$ws.Range("E6").Value = "  +2.41%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "0.527"
$ws.Range("E8").Value = "  +1.20%  "
$ws.Range("D9").Value = "2.626.06"
$ws.Range("E9").Value = "  +3.97%  "
$ws.Range("D10").Value = "0.167"
$ws.Range("E10").Value = "  +13.92%  "
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("D12").Value = "0.348"
$ws.Range("E12").Value = "  +2.42%  "
$ws.Range("D13").Value = "5.02"
$ws.Range("E13").Value = "  +0.55%  "
$ws.Range("D14").Value = "3.114.18"
$ws.Range("E14").Value = "  +4.20%  "
$ws.Range("D15").Value = "0.0000188"
$ws.Range("E15").Value = "  +10.13%  "
$ws.Range("D16").Value = "26.77"
$ws.Range("E16").Value = "  +2.00%  "
$ws.Range("D17").Value = "71.369.67"
$ws.Range("E17").Value = "  +4.07%  "
$ws.Range("D18").Value = "2.625.44"
$ws.Range("E18").Value = "  +4.28%  "
$ws.Range("D19").Value = "382.12"
$ws.Range("E19").Value = "  +8.33%  "
$ws.Range("D20").Value = "7.92"
$ws.Range("E20").Value = "  +5.44%  "
$ws.Range("E21").Value = "  +3.71%  "
$ws.Range("D22").Value = "4.15"
$ws.Range("E22").Value = "  -1.44%  "
$ws.Range("D23").Value = "72.74"
$ws.Range("E23").Value = "  +2.64%  "
$ws.Range("B24").Value = "SuiNetwork"
$ws.Range("C24").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D24").Value = "1.98"
$ws.Range("E24").Value = "  +16.40%  "
$ws.Range("B25").Value = "NEARProtocol"
$ws.Range("C25").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D25").Value = "4.48"
$ws.Range("E25").Value = "  +5.87%  "
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("D27").Value = "9.75"
$ws.Range("E27").Value = "  +8.38%  "
$ws.Range("D28").Value = "2.761.80"
$ws.Range("E28").Value = "  +2.64%  "
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.44%  "
$ws.Range("E30").Value = "  +8.17%  "
$ws.Range("D31").Value = "550.91"
$ws.Range("E31").Value = "  +8.38%  "
$ws.Range("D32").Value = "8.08"
$ws.Range("E32").Value = "  +3.40%  "
$ws.Range("D33").Value = "1.33"
$ws.Range("E33").Value = "  +5.97%  "
$ws.Range("E34").Value = "  +3.43%  "
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("E36").Value = "  +1.93%  "
$ws.Range("D37").Value = "0.117"
$ws.Range("E37").Value = "  -2.50%  "
$ws.Range("D38").Value = "19.29"
$ws.Range("E38").Value = "  +4.78%  "
$ws.Range("D39").Value = "1.90"
$ws.Range("E39").Value = "  +7.24%  "
$ws.Range("D40").Value = "19.05"
$ws.Range("E40").Value = "  +1.95%  "
$ws.Range("D41").Value = "1.39"
$ws.Range("E41").Value = "  +4.97%  "
$ws.Range("D42").Value = "2.65"
$ws.Range("E42").Value = "  +9.34%  "
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("D44").Value = "5.06"
$ws.Range("E44").Value = "  +4.37%  "
$ws.Range("D45").Value = "0.332"
$ws.Range("E45").Value = "  +1.87%  "
$ws.Range("D46").Value = "40.02"
$ws.Range("E46").Value = "  +2.55%  "
$ws.Range("D47").Value = "154.41"
$ws.Range("E47").Value = "  +1.07%  "
$ws.Range("D48").Value = "3.65"
$ws.Range("E48").Value = "  +2.16%  "
$ws.Range("E50").Value = "  +2.64%  "
$ws.Range("E51").Value = "  +1.72%  "
